$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.051.91'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '3.064.04'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '560.31'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').Value = '143.57'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '3.061.59'
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E9').Value = '  +3.42%  '
$ws.Range('D10').Value = '0.155'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').Value = '6.16'
$ws.Range('E11').Value = '  -10.19%  '
$ws.Range('D12').Value = '0.497'
$ws.Range('E12').Value = '  +8.67%  '
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('D14').Value = '35.83'
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').Value = '3.565.82'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '64.095.74'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '3.067.24'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('D20').Value = '478.65'
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').Value = '13.97'
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('E22').Value = '  +1.57%  '
$ws.Range('D23').Value = '14.38'
$ws.Range('E23').Value = '  +8.79%  '
$ws.Range('D24').Value = '7.57'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').Value = '82.69'
$ws.Range('E25').Value = '  +2.02%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '2.81'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').Value = '8.11'
$ws.Range('E28').Value = '  +2.02%  '
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('D31').Value = '26.33'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').Value = '5.81'
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('D35').Value = '6.25'
$ws.Range('E35').Value = '  +2.63%  '
$ws.Range('D36').Value = '54.84'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.0412'
$ws.Range('E37').Value = '  +1.57%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '452.46'
$ws.Range('E38').Value = '  -2.21%  '
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('E40').Value = '  +3.55%  '
$ws.Range('D41').Value = '3.029.69'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '8.29'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('E44').Value = '  +2.87%  '
$ws.Range('D45').Value = '27.74'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  +6.59%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  +1.86%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '119.09'
$ws.Range('E49').Value = '  +1.82%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0518'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('E51').Value = '  +0.45%  '
